# Refresh the "cryptos" price/volume snapshot (rows 2-51) to the values
# from the latest GitHub Actions run. Price (col D) and Volume(1h) (col E)
# are plain text in this sheet (not real numbers), so for any new price
# that happens to look numeric ("1.02", "0.991", "0.0000127", ...) we
# briefly force Text number-formatting before the assignment and restore
# the cell's original style afterwards - otherwise Excel's COM Value
# setter would silently coerce the string to a Double and mangle things
# like trailing zeros / leading zeros / thousand-dot-separated prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '53.959.77'
$ws.Range("E2").Value = '  -9.49%  '
$ws.Range("D3").Value = '2.394.25'
$ws.Range("E3").Value = '  -16.42%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue "D5" '459.27'
$ws.Range("E5").Value = '  -9.88%  '
Set-TextValue "D6" '129.58'
$ws.Range("E6").Value = '  -6.38%  '
Set-TextValue "D7" '0.995'
$ws.Range("E7").Value = '  -0.60%  '
Set-TextValue "D8" '0.484'
$ws.Range("E8").Value = '  -8.76%  '
$ws.Range("D9").Value = '2.400.87'
$ws.Range("E9").Value = '  -16.03%  '
Set-TextValue "D10" '0.0937'
$ws.Range("E10").Value = '  -9.57%  '
Set-TextValue "D11" '5.34'
$ws.Range("E11").Value = '  -11.85%  '
Set-TextValue "D12" '0.317'
$ws.Range("E12").Value = '  -9.20%  '
$ws.Range("E13").Value = '  -4.35%  '
$ws.Range("D14").Value = '2.800.71'
$ws.Range("E14").Value = '  -16.72%  '
$ws.Range("D15").Value = '53.803.29'
$ws.Range("E15").Value = '  -10.25%  '
Set-TextValue "D16" '19.51'
$ws.Range("E16").Value = '  -10.93%  '
Set-TextValue "D17" '0.0000127'
$ws.Range("E17").Value = '  -7.28%  '
$ws.Range("D18").Value = '2.397.80'
$ws.Range("E18").Value = '  -16.04%  '
Set-TextValue "D19" '4.20'
$ws.Range("E19").Value = '  -12.01%  '
Set-TextValue "D20" '310.54'
$ws.Range("E20").Value = '  -11.44%  '
Set-TextValue "D21" '9.25'
$ws.Range("E21").Value = '  -17.44%  '
Set-TextValue "D22" '1.02'
$ws.Range("E22").Value = '  +1.81%  '
Set-TextValue "D23" '5.66'
$ws.Range("E23").Value = '  +0.29%  '
Set-TextValue "D24" '5.38'
$ws.Range("E24").Value = '  -15.07%  '
Set-TextValue "D25" '55.88'
$ws.Range("E25").Value = '  -11.68%  '
Set-TextValue "D26" '1.01'
$ws.Range("E26").Value = '  +0.63%  '
Set-TextValue "D27" '0.378'
$ws.Range("E27").Value = '  -13.24%  '
Set-TextValue "D28" '0.152'
$ws.Range("E28").Value = '  -12.61%  '
$ws.Range("D29").Value = '2.464.79'
$ws.Range("E29").Value = '  -17.72%  '
Set-TextValue "D30" '7.11'
$ws.Range("E30").Value = '  -5.98%  '
Set-TextValue "D31" '0.996'
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("D32").Value = '0.0₃0708'
$ws.Range("E32").Value = '  -14.97%  '
Set-TextValue "D33" '147.76'
$ws.Range("E33").Value = '  -2.65%  '
Set-TextValue "D34" '17.61'
$ws.Range("E34").Value = '  -7.88%  '
Set-TextValue "D35" '1.39'
$ws.Range("E35").Value = '  -14.91%  '
Set-TextValue "D36" '4.98'
$ws.Range("E36").Value = '  -7.90%  '
Set-TextValue "D37" '3.46'
$ws.Range("E37").Value = '  -18.29%  '
Set-TextValue "D38" '1.04'
$ws.Range("E38").Value = '  -10.77%  '
Set-TextValue "D39" '0.789'
$ws.Range("E39").Value = '  -17.33%  '
Set-TextValue "D40" '33.48'
$ws.Range("E40").Value = '  -9.42%  '
Set-TextValue "D41" '0.991'
$ws.Range("E41").Value = '  -0.61%  '
Set-TextValue "D42" '0.598'
$ws.Range("E42").Value = '  -6.13%  '
Set-TextValue "D43" '3.28'
$ws.Range("E43").Value = '  -7.72%  '
Set-TextValue "D44" '0.0525'
$ws.Range("E44").Value = '  -7.18%  '
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("D46").Value = '1.974.24'
$ws.Range("E46").Value = '  -11.81%  '
Set-TextValue "D47" '1.22'
$ws.Range("E47").Value = '  -13.42%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D48" '0.0215'
$ws.Range("E48").Value = '  -5.92%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D49" '0.0856'
$ws.Range("E49").Value = '  -4.20%  '
Set-TextValue "D50" '4.29'
$ws.Range("E50").Value = '  -8.13%  '
$ws.Range("E51").Value = '  -17.35%  '
